$d = $word.ActiveDocument

# --- 1. Add new paragraph style "Abstract Title" (styleId AbstractTitle),
#        based on Normal, next paragraph style Abstract, inserted ahead of
#        the "Abstract" style in intent (styles collection appends, but the
#        formatting/linking below reproduces the authored style exactly).
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles.Item("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles.Item("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 2. "Abstract" style: reduce space-before from 15pt (300 twips) to
#        5pt (100 twips); space-after stays 15pt (300 twips).
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. "ImportTok" character style: add green + bold run formatting.
$importTok = $d.Styles.Item("ImportTok")
$importTok.Font.Color = 32768
$importTok.Font.Bold = $true

# --- 4. "BuiltInTok" character style: add green run formatting.
$builtInTok = $d.Styles.Item("BuiltInTok")
$builtInTok.Font.Color = 32768
